$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A34:E34").NumberFormat = "@"

$ws.Range("A34").Value = "2025-10-20"
$ws.Range("B34").Value = "Pick 3"
$ws.Range("C34").Value = "251020"
$ws.Range("D34").Value = "9-4-6"
$ws.Range("E34").Value = "2025-10-20T21:38:22.367+04:00"
